$d = $word.ActiveDocument

# Replace text of paragraphs 2 through 19 (content paragraphs 1-18)
$d.Paragraphs.Item(2).Range.Text = 'If you see a severity one or two, take a screenshot of it and post it inside the AMS standup group or e-services group. And then you can address to UN, Eugene or the team leads. Onboarding of boarding, there is a main folder onboarding of boarding inside that KT to Eugene and UN.'
$d.Paragraphs.Item(3).Range.Text = 'Once it gets stabilized, then then yeah, you can able to manage not an issue. Once you are familiar you can take this opportunity to onboard them. Meanwhile when you are onboarding if you find any difficulties you message me. So after that we can keep the I mean make it completed. So the next item shall we go to the resource talk take.'
$d.Paragraphs.Item(4).Range.Text = 'This is for DXE. Onshore people. And snow access. I think you got snow access but I need to confirm whether you can able to log in. So let''s make it as like blank. You can add them later. Do you know where are you from in India or the south or the north? No, Singapore. Oh, my place is it? Yes. Are you in no place? Oh, I''m from southern part.'
$d.Paragraphs.Item(5).Range.Text = 'The severity level is one, two or three, depending on the severity of the issue. If it is a complete application down and none of the users unable to log in, that is considered a critical issue. It is also considered a security related issue if it is due to some security or data breach. If the use shoe is brought up by user, then it has to be like, you know, the, the reported by must be username.'
$d.Paragraphs.Item(6).Range.Text = 'NTUC will be connecting to, to NTU CTO and who is the IT team. They will coordinate with BTB, the infra team and DXC application team. Once we get a severity one, we, the snow ticket has to be marked as severity one. If it is a severity two, actually they will do, the helpers will do but sometimes they will forget.'
$d.Paragraphs.Item(7).Range.Text = 'If you are seeing anything like something not usual, something that, like alarming, which like multiple users impacted, then you need to take a screenshot and put it in a MS group, stand-up group and address to Eugene. So then based on their confirmation, you can, later you can check whether the ticket is categorized correctly. Then within 30 minutes, right? Oh, sorry, within one hour, we need to start sending this template.'
$d.Paragraphs.Item(8).Range.Text = 'Every one hour we need to send and finally issue resolved. But we get very, very less severity one and two, very, like hardly we get. So it''s not like often. So yeah, but this is the process of severity 1 and two. Okay. Very clear. Any questions? No questions so far. I want one question. This is, Eugene is from technical team. Eugene is the technical lead, MS technical lead. Under him, there are four technical leads, like I will share you their names later.'
$d.Paragraphs.Item(9).Range.Text = 'You can use this highlighter also. Then put it inside the MS. If you think like they are not responding to code, where is the MS? Okay. So the monitor complete tracking SLA and ensure it meets the month end. This is the my NTUC CR team. And this is E-Service and U-Tab lead. And he is the SharePoint. Deepika is CR. Kong is the. my NT UC team did. UN is overall. And TIN is e-Services.'
$d.Paragraphs.Item(10).Range.Text = 'This is the one of the update, status update. Once we get the resolution time, we need to put in in the same format like the check the. resolution time with the technical team before like 20, maybe, or 10, 11, something, then other comments, the issue resolved. Here, the issues resolved after rebooting the server or something like that. We need to. put and then if we are providing the RCA, you check, better you check with Eugene or the technical RCA will be provided late.'
$d.Paragraphs.Item(11).Range.Text = 'There are two weekly meetings, okay. One is the AMS meeting. There we will provide all the maintenance service related things. And another is the CR meeting. In that meeting, we will discuss all the PMs will join, also the PM''s head will join. From DXC side, DXC management will join Bala. Next to Bala is Sadat, sometimes he will also join. And Deepika will provide updates.'
$d.Paragraphs.Item(12).Range.Text = 'So, last week, you see here, there is a follow-up action item. In the meeting, NTUCTO, Josephine suggested to keep this severity one, follow- up action is in this weekly slide until we close all the action items. So, this slide, you no need to worry. This is updated by the technical team members. But in case, if you are asked to present this slide,. you should know what is inside.'
$d.Paragraphs.Item(13).Range.Text = 'Each incident, for example, if we had two or three incident last week, then we had to put each one. provide. Usually for one month, we have, like, you see here, in November. And after that, this is in December, only one incident, and January one. So, every month one and sometimes last from July to November, we didn''t even have any incidents.'
$d.Paragraphs.Item(14).Range.Text = 'If you see such slides, or this one reduction, this one suddenly again increase, just understand if you ask the mind to see team lead, they can able to tell you. Again, the team leads will also join the meeting. So if you are stuck, they will answer, no need to worry about it. So out of this 51, 49 and close, resolved.'
$d.Paragraphs.Item(15).Range.Text = 'Less than 30 days as of 26 February, less than 30 day ticket, okay, this is the trend. No aging ticket more than. the one, only one awaiting customer. So this is related with the mind you see, UAT bug fix. So we are waiting for customer to sign off. All these will be provided by the support team. You just need to ensure like you are.'
$d.Paragraphs.Item(16).Range.Text = 'You need to make sure these are like, you know, if you are presenting, you should be knowing this. If there is no update from January of February, then you need to ask them, hey, can you please fill up the latest update until today, what''s happened? And yeah, this is also the same, all these tickets, they have to provide. And this is the service improvement.'
$d.Paragraphs.Item(17).Range.Text = ' NGUC asks us to do some enhancement to the current environment, to our current application, utilizing the five mandates. as service improvement to NGUC, okay? So for a case upgrade, okay, maybe sometimes nine mandates, but they will minus five, minus nine, four, four in the next month. So if this is accumulated, like five times 12, it''s accumulated. So the total available mandates must be this.'
$d.Paragraphs.Item(18).Range.Text = 'The vendor site, customer site. And here, this is the monthly, okay. Let me show. MS and CR weekly update. And this is, okay, number one, age ticket details for this one. More than 30 days. If they ask something about the ticket details and if they want DXC to do some follow-up action item or something to take note, then we need to put it inside.'
$d.Paragraphs.Item(19).Range.Text = 'When we repair the meeting minutes, then we had to check, double check with them. And then aging ticket process confirmation by in the meeting, it was mentioned that the older tickets, DXC has completed the task, fix the bug, but it is awaiting confirmation from the user. Victor proposed that the DXC should close this ticket if there is no response from the business user after several follow-ups.'

# Replace text of the last paragraph (29)
$d.Paragraphs.Item(29).Range.Text = 'Every week, three times, we have a call. Whenever we have new updates, we will internally review an update on Wednesday, because Thursday is the meeting with NTUC. On Wednesday, we stop updating and we will send the updated slide and then we will go through the, that slides in the Thursday meeting.'

# Delete paragraphs 20 through 28 (9 paragraphs that are removed entirely)
# Delete from highest index to lowest to keep indices stable
for ($i = 28; $i -ge 20; $i--) {
    $d.Paragraphs.Item($i).Range.Delete()
}
